$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 80
$ws.Range("H80").Value = 314.70587
$ws.Range("I80").Value = 300.0909
$ws.Range("J80").Value = 341.5
$ws.Range("K80").Value = 900.2727
$ws.Range("L80").Value = 1024.5
$ws.Range("M80").Value = 97.72730000000001
$ws.Range("N80").Value = -3020.5

# Row 83
$ws.Range("H83").Value = 314.70587
$ws.Range("I83").Value = 300.0909
$ws.Range("J83").Value = 341.5
$ws.Range("K83").Value = 2700.8181
$ws.Range("L83").Value = 3073.5
$ws.Range("M83").Value = 2291.1819
$ws.Range("N83").Value = -13057.5

# Row 88
$ws.Range("H88").Value = 3033.3333
$ws.Range("I88").Value = 1900
$ws.Range("J88").Value = 3100
$ws.Range("K88").Value = 1900
$ws.Range("L88").Value = 3100
$ws.Range("M88").Value = -1494
$ws.Range("N88").Value = -3912

# Row 91
$ws.Range("H91").Value = 3033.3333
$ws.Range("I91").Value = 1900
$ws.Range("J91").Value = 3100
$ws.Range("K91").Value = 1900
$ws.Range("L91").Value = 3100
$ws.Range("M91").Value = -496
$ws.Range("N91").Value = -5908


$ws = $wb.Worksheets.Item("ARM")
# Row 88
$ws.Range("H88").Value = 3552.8333
$ws.Range("I88").Value = 2663.4
$ws.Range("K88").Value = 2663.4
$ws.Range("M88").Value = -2257.4

# Row 91
$ws.Range("H91").Value = 3552.8333
$ws.Range("I91").Value = 2663.4
$ws.Range("K91").Value = 2663.4
$ws.Range("M91").Value = -1259.4


$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value = 1962.3928
$ws.Range("I86").Value = 1723.8695
$ws.Range("J86").Value = 3059.6
$ws.Range("K86").Value = 1723.8695
$ws.Range("L86").Value = 3059.6
$ws.Range("M86").Value = -600.8695
$ws.Range("N86").Value = -5305.6

# Row 89
$ws.Range("H89").Value = 1962.3928
$ws.Range("I89").Value = 1723.8695
$ws.Range("J89").Value = 3059.6
$ws.Range("K89").Value = 8619.3475
$ws.Range("L89").Value = 15298
$ws.Range("M89").Value = -3003.3475
$ws.Range("N89").Value = -26530

# Row 131
$ws.Range("H131").Value = 43852
$ws.Range("J131").Value = 43852
$ws.Range("L131").Value = 43852
$ws.Range("N131").Value = -53932


$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 1335.909
$ws.Range("I31").Value = 1113.8096
$ws.Range("J31").Value = 6000
$ws.Range("K31").Value = 1113.8096
$ws.Range("L31").Value = 6000
$ws.Range("M31").Value = -818.8096
$ws.Range("N31").Value = -6590

# Row 34
$ws.Range("H34").Value = 1335.909
$ws.Range("I34").Value = 1113.8096
$ws.Range("J34").Value = 6000
$ws.Range("K34").Value = 1113.8096
$ws.Range("L34").Value = 6000
$ws.Range("M34").Value = -911.8096
$ws.Range("N34").Value = -6404

# Row 62
$ws.Range("H62").Value = 3402.5
$ws.Range("I62").Value = 3120
$ws.Range("J62").Value = 4250
$ws.Range("K62").Value = 3120
$ws.Range("L62").Value = 4250
$ws.Range("M62").Value = -2496
$ws.Range("N62").Value = -5498

# Row 65
$ws.Range("H65").Value = 3402.5
$ws.Range("I65").Value = 3120
$ws.Range("J65").Value = 4250
$ws.Range("K65").Value = 15600
$ws.Range("L65").Value = 21250
$ws.Range("M65").Value = -12480
$ws.Range("N65").Value = -27490


$ws = $wb.Worksheets.Item("CUL")
# Row 86
$ws.Range("H86").Value = 200
$ws.Range("I86").Value = 200
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 600
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = $null
$ws.Range("N86").Value = 586

# Row 89
$ws.Range("H89").Value = 200
$ws.Range("I89").Value = 200
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 1800
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = $null
$ws.Range("N89").Value = 4128

# Row 131
$ws.Range("H131").Value = 1431010
$ws.Range("J131").Value = 2002115.8
$ws.Range("L131").Value = 6006347.4
$ws.Range("N131").Value = -6016427.4


$ws = $wb.Worksheets.Item("GSM")
# Row 51
$ws.Range("H51").Value = 44496.332
$ws.Range("I51").Value = 41000
$ws.Range("J51").Value = 45195.6
$ws.Range("K51").Value = 41000
$ws.Range("L51").Value = 45195.6
$ws.Range("M51").Value = -40491
$ws.Range("N51").Value = -46213.6

# Row 113
$ws.Range("H113").Value = 1673.0667
$ws.Range("J113").Value = 950
$ws.Range("L113").Value = 950
$ws.Range("N113").Value = -5290


$ws = $wb.Worksheets.Item("LTW")
# Row 62
$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = $null
$ws.Range("M62").Value = $null
$ws.Range("N62").Value = 0

# Row 65
$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = $null
$ws.Range("M65").Value = $null
$ws.Range("N65").Value = 0

# Row 68
$ws.Range("H68").Value = 1531.3572
$ws.Range("I68").Value = 1405.4546
$ws.Range("J68").Value = 1993
$ws.Range("K68").Value = 1405.4546
$ws.Range("L68").Value = 1993
$ws.Range("M68").Value = -656.4546
$ws.Range("N68").Value = -3491

# Row 71
$ws.Range("H71").Value = 1531.3572
$ws.Range("I71").Value = 1405.4546
$ws.Range("J71").Value = 1993
$ws.Range("K71").Value = 7027.273
$ws.Range("L71").Value = 9965
$ws.Range("M71").Value = -3283.273
$ws.Range("N71").Value = -17453


$ws = $wb.Worksheets.Item("WVR")
# Row 62
$ws.Range("H62").Value = 5362.5
$ws.Range("I62").Value = 3600
$ws.Range("J62").Value = 7125
$ws.Range("K62").Value = 3600
$ws.Range("L62").Value = 7125
$ws.Range("M62").Value = -2976
$ws.Range("N62").Value = -8373

# Row 65
$ws.Range("H65").Value = 5362.5
$ws.Range("I65").Value = 3600
$ws.Range("J65").Value = 7125
$ws.Range("K65").Value = 18000
$ws.Range("L65").Value = 35625
$ws.Range("M65").Value = -14880
$ws.Range("N65").Value = -41865

# Row 81
$ws.Range("H81").Value = 4538.4
$ws.Range("I81").Value = 1975
$ws.Range("J81").Value = 5470.5454
$ws.Range("K81").Value = 3950
$ws.Range("L81").Value = 10941.0908
$ws.Range("M81").Value = -2889
$ws.Range("N81").Value = -13063.0908

# Row 84
$ws.Range("H84").Value = 4538.4
$ws.Range("I84").Value = 1975
$ws.Range("J84").Value = 5470.5454
$ws.Range("K84").Value = 19750
$ws.Range("L84").Value = 54705.454
$ws.Range("M84").Value = -14446
$ws.Range("N84").Value = -65313.454

